$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2032520325203252
$ws.Range("C2").Value = 0.5447154471544715
$ws.Range("P2").Value = 0.1626016260162602
$ws.Range("S2").Value = 0.08943089430894309
$ws.Range("C3").Value = 0.03472222222222222
$ws.Range("J3").Value = 0.02777777777777778
$ws.Range("P3").Value = 0.75
$ws.Range("S3").Value = 0.1875
$ws.Range("J4").Value = 0.03773584905660377
$ws.Range("P4").Value = 0.6415094339622641
$ws.Range("S4").Value = 0.3207547169811321
$ws.Range("B6").Value = 0.03930131004366812
$ws.Range("D6").Value = 0.02183406113537118
$ws.Range("F6").Value = 0.07423580786026202
$ws.Range("J6").Value = 0.2489082969432314
$ws.Range("O6").Value = 0.01746724890829694
$ws.Range("Q6").Value = 0.1615720524017467
$ws.Range("R6").Value = 0.07423580786026202
$ws.Range("S6").Value = 0.3624454148471616
$ws.Range("B7").Value = 0.07731958762886598
$ws.Range("D7").Value = 0.02577319587628866
$ws.Range("F7").Value = 0.07216494845360824
$ws.Range("J7").Value = 0.1030927835051546
$ws.Range("O7").Value = 0.02577319587628866
$ws.Range("Q7").Value = 0.154639175257732
$ws.Range("R7").Value = 0.09793814432989691
$ws.Range("S7").Value = 0.4432989690721649
$ws.Range("B8").Value = 0.08433734939759036
$ws.Range("D8").Value = 0.02208835341365462
$ws.Range("E8").Value = 0.002008032128514056
$ws.Range("F8").Value = 0.04819277108433735
$ws.Range("J8").Value = 0.1104417670682731
$ws.Range("O8").Value = 0.01004016064257028
$ws.Range("Q8").Value = 0.1867469879518072
$ws.Range("R8").Value = 0.09437751004016064
$ws.Range("S8").Value = 0.4417670682730924
$ws.Range("B9").Value = 0.0707070707070707
$ws.Range("D9").Value = 0.01515151515151515
$ws.Range("F9").Value = 0.03535353535353535
$ws.Range("J9").Value = 0.101010101010101
$ws.Range("O9").Value = 0.0101010101010101
$ws.Range("Q9").Value = 0.1767676767676768
$ws.Range("R9").Value = 0.101010101010101
$ws.Range("S9").Value = 0.4898989898989899
$ws.Range("B10").Value = 0.08702290076335878
$ws.Range("D10").Value = 0.02519083969465649
$ws.Range("E10").Value = 0.0007633587786259542
$ws.Range("F10").Value = 0.0633587786259542
$ws.Range("J10").Value = 0.1206106870229008
$ws.Range("O10").Value = 0.01374045801526718
$ws.Range("Q10").Value = 0.2053435114503817
$ws.Range("R10").Value = 0.1015267175572519
$ws.Range("S10").Value = 0.3824427480916031
$ws.Range("G11").Value = 0.1563517915309446
$ws.Range("J11").Value = 0.08143322475570032
$ws.Range("K11").Value = 0.2149837133550489
$ws.Range("L11").Value = 0.5276872964169381
$ws.Range("S11").Value = 0.01954397394136808
$ws.Range("G12").Value = 0.7349397590361446
$ws.Range("J12").Value = 0.1987951807228916
$ws.Range("K12").Value = 0.01807228915662651
$ws.Range("L12").Value = 0.01204819277108434
$ws.Range("S12").Value = 0.03614457831325301
$ws.Range("F15").Value = 0.02272727272727273
$ws.Range("H15").Value = 0.2409090909090909
$ws.Range("I15").Value = 0.05909090909090909
$ws.Range("J15").Value = 0.3227272727272728
$ws.Range("K15").Value = 0.07272727272727272
$ws.Range("O15").Value = 0.08636363636363636
$ws.Range("S15").Value = 0.1954545454545455
$ws.Range("F16").Value = 0.01704545454545454
$ws.Range("H16").Value = 0.1704545454545454
$ws.Range("I16").Value = 0.06818181818181818
$ws.Range("J16").Value = 0.4318181818181818
$ws.Range("K16").Value = 0.08522727272727272
$ws.Range("M16").Value = 0.02272727272727273
$ws.Range("O16").Value = 0.07386363636363637
$ws.Range("S16").Value = 0.1306818181818182
$ws.Range("F17").Value = 0.02183406113537118
$ws.Range("H17").Value = 0.1986899563318777
$ws.Range("I17").Value = 0.1048034934497817
$ws.Range("J17").Value = 0.4104803493449782
$ws.Range("K17").Value = 0.09388646288209607
$ws.Range("M17").Value = 0.01746724890829694
$ws.Range("N17").Value = 0.002183406113537118
$ws.Range("O17").Value = 0.05458515283842795
$ws.Range("S17").Value = 0.09606986899563319
$ws.Range("F18").Value = 0.0211864406779661
$ws.Range("H18").Value = 0.1694915254237288
$ws.Range("I18").Value = 0.1059322033898305
$ws.Range("J18").Value = 0.3771186440677966
$ws.Range("K18").Value = 0.1101694915254237
$ws.Range("M18").Value = 0.01694915254237288
$ws.Range("O18").Value = 0.07203389830508475
$ws.Range("S18").Value = 0.1271186440677966
$ws.Range("F19").Value = 0.02296296296296296
$ws.Range("H19").Value = 0.2133333333333333
$ws.Range("I19").Value = 0.0725925925925926
$ws.Range("J19").Value = 0.38
$ws.Range("K19").Value = 0.1007407407407407
$ws.Range("M19").Value = 0.02148148148148148
$ws.Range("O19").Value = 0.06592592592592593
$ws.Range("S19").Value = 0.122962962962963
